$d = $word.ActiveDocument

# 1. Remove the "Convert standalone iterators..." bullet entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Convert standalone iterators to be created by factory class*") {
        $p.Range.Delete()
        break
    }
}

# 2. Add a new "Light Windows API wrapper." bullet right after the
#    "Stack trace in exception object." bullet (before the "New Modules" heading).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Stack trace in exception object*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Light Windows API wrapper."
        break
    }
}
